$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) and "全部类型" (All Types) contain the same
# event-list data; update the "想去人数" (wish-to-go count) figures for
# the rows that changed in the latest data pull.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 1737
    $ws.Range("F5").Value = 6266
    $ws.Range("F6").Value = 102
}
